# Applies the coinranking "cryptos" price/volume refresh for this run.
# Columns B/C (coin name, link) are plain text and can be set directly.
# Columns D/E (price, % volume) are numeric-looking strings (e.g. "1.00",
# "8.00", thousand-dot-grouped prices, "  +1.29%  ") that Excel's COM
# auto-detection would otherwise silently coerce into real numbers/percents
# (losing the exact text, e.g. "1.00" -> 1). Forcing the cell to Text via
# NumberFormat "@" before the write preserves the literal string, and
# resetting the cell Style back to "Normal" afterwards drops the temporary
# text-format styling so the cell ends up unstyled again, same as before.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '67.930.14'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.517.69'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '600.89'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +1.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '181.45'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +4.75%  '
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.518.51'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c.Style = 'Normal'

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.596'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.141'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +7.32%  '
$c.Style = 'Normal'

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '7.16'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -1.78%  '
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +0.80%  '
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.122.72'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +0.18%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '32.72'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +12.46%  '
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '67.881.41'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +1.23%  '
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +0.53%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '3.515.01'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.81%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '6.36'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +0.52%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '14.49'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +2.36%  '
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '400.38'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +1.37%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '8.00'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -0.50%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '73.71'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.76%  '
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.40%  '
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +1.60%  '
$c.Style = 'Normal'

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '10.56'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +2.75%  '
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -2.59%  '
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '
$c.Style = 'Normal'

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '6.27'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -0.59%  '
$c.Style = 'Normal'

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.46'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +0.23%  '
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +1.56%  '
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '23.92'
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.41%  '
$c.Style = 'Normal'

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '7.51'
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +1.70%  '
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -3.06%  '
$c.Style = 'Normal'

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '163.30'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -0.16%  '
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -0.20%  '
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +0.82%  '
$c.Style = 'Normal'

$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.78'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +7.69%  '
$c.Style = 'Normal'

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '7.01'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.71%  '
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.890.40'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +2.75%  '
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '4.70'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +0.28%  '
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0736'
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -1.58%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '26.34'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -0.31%  '
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '26.88'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -2.32%  '
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '42.35'
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -0.83%  '
$c.Style = 'Normal'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.0304'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '345.38'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +1.90%  '
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.19%  '
$c.Style = 'Normal'
